# "Added MySQL dynamic query for variable rows number"
# Duplicate the existing A:C data block into a new D:F block on the same
# rows, after shifting the original data down by one row (rows 1-12 -> 2-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the existing data down one row (row 1 -> row 2, ... row 12 -> row 13)
$ws.Rows.Item(1).Insert()

# Mirror the A2:C13 block into D2:F13 so each row carries two copies of the
# record (used to drive a MySQL query across a variable number of rows).
$src = $ws.Range("A2:C13")
$dst = $ws.Range("D2:F13")
$src.Copy($dst)

# Column E (the copied date column) should display with the same width as
# column B.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Leave the sheet with the whole new block selected (cursor lands on the
# first/top-left cell of the selection, A2).
$ws.Range("A2:F13").Select()
